# Apply the edit described by the diff:
#  - Remove the "Other" row (original row 4: Other / 351.588884563236 / 0.814919535887345)
#    which shifts the subsequent rows (Production areas, Protected areas,
#    Requires individual assessment) up by one.
#  - Update the values of the new last row (originally "Requires individual
#    assessment" row) to the new B/C values from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 4 ("Other"), shifting rows 5-7 up to 4-6.
$ws.Rows.Item(4).Delete()

# After the shift, row 6 is now "Requires individual assessment" and its
# numeric values need to be updated to the new figures from the diff.
$ws.Range("B6").Value = 2529.90342608951
$ws.Range("C6").Value = 5.86385922976431
